# MIAME-ENV_data.xlsx: apply new naming rules (spaces -> underscores) to
# sheet tab names and to the entity/attribute metadata rows that mirror
# those names (commit: "updated MIAME-ENV for new naming rules").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "entities" sheet: entity `name` column (A) loses its spaces.
# ---------------------------------------------------------------------
$wsEntities = $wb.Worksheets.Item("entities")
$wsEntities.Range("A2").Value = "General_features"
$wsEntities.Range("A3").Value = "Contact_Person"
$wsEntities.Range("A4").Value = "Key_Concepts"
$wsEntities.Range("A7").Value = "Location_of_Documents"
# A5 "Concept" and A6 "Bibliography" are unchanged.

# ---------------------------------------------------------------------
# 2) "attributes" sheet: attribute `name` column (A) loses its spaces
#    (and apostrophes), and the `entity` column (B) is updated to the
#    renamed entity ids.
# ---------------------------------------------------------------------
$wsAttributes = $wb.Worksheets.Item("attributes")

# column A (attribute name) renames
$wsAttributes.Range("A4").Value  = "Document_Type"
$wsAttributes.Range("A5").Value  = "Group_"
$wsAttributes.Range("A6").Value  = "Main_Website"
$wsAttributes.Range("A7").Value  = "MI_Checklist_s_Name"
$wsAttributes.Range("A8").Value  = "MI_Checklist_s_Acronym"
$wsAttributes.Range("A9").Value  = "Current_Version_Designation"
$wsAttributes.Range("A10").Value = "ReleaseDate_Current_Version"
$wsAttributes.Range("A11").Value = "General_Comments"
$wsAttributes.Range("A12").Value = "Full_Name"
$wsAttributes.Range("A13").Value = "Email_Address"
$wsAttributes.Range("A18").Value = "PubMed_Identifier"
$wsAttributes.Range("A19").Value = "Digital_Object_Identifier"

# column B (owning entity) renames - mirrors the entity id renames above
$wsAttributes.Range("B2").Value  = "MIAMEENV_General_features"
$wsAttributes.Range("B3").Value  = "MIAMEENV_General_features"
$wsAttributes.Range("B4").Value  = "MIAMEENV_General_features"
$wsAttributes.Range("B5").Value  = "MIAMEENV_General_features"
$wsAttributes.Range("B6").Value  = "MIAMEENV_General_features"
$wsAttributes.Range("B7").Value  = "MIAMEENV_General_features"
$wsAttributes.Range("B8").Value  = "MIAMEENV_General_features"
$wsAttributes.Range("B9").Value  = "MIAMEENV_General_features"
$wsAttributes.Range("B10").Value = "MIAMEENV_General_features"
$wsAttributes.Range("B11").Value = "MIAMEENV_General_features"

$wsAttributes.Range("B12").Value = "MIAMEENV_Contact_Person"
$wsAttributes.Range("B13").Value = "MIAMEENV_Contact_Person"

$wsAttributes.Range("B14").Value = "MIAMEENV_Key_Concepts"
$wsAttributes.Range("B15").Value = "MIAMEENV_Key_Concepts"

$wsAttributes.Range("B21").Value = "MIAMEENV_Location_of_Documents"
$wsAttributes.Range("B22").Value = "MIAMEENV_Location_of_Documents"
$wsAttributes.Range("B23").Value = "MIAMEENV_Location_of_Documents"

# ---------------------------------------------------------------------
# 3) Rename the worksheet tabs themselves (do this last so the
#    name-based lookups above keep working).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("MIAMEENV_General features").Name = "MIAMEENV_General_features"
$wb.Worksheets.Item("MIAMEENV_Contact Person").Name = "MIAMEENV_Contact_Person"
$wb.Worksheets.Item("MIAMEENV_Key Concepts").Name = "MIAMEENV_Key_Concepts"
$wb.Worksheets.Item("MIAMEENV_Location of Documents").Name = "MIAMEENV_Location_of_Documents"

# ---------------------------------------------------------------------
# 4) Move the active tab/selection to the (now first) general-features
#    sheet, matching the post-edit view state.
# ---------------------------------------------------------------------
$wsGeneral = $wb.Worksheets.Item("MIAMEENV_General_features")
$wsGeneral.Activate() | Out-Null
$wsGeneral.Range("I2").Select() | Out-Null
